$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.101371
$ws.Range("H2").Value = 0.304113
$ws.Range("I2").Value = 0.0004873780190420389
$ws.Range("J2").Value = 0.0004888126654476159
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.343285
$ws.Range("N2").Value = 1.029855
$ws.Range("O2").Value = 0.9867450936054706
$ws.Range("P2").Value = 0.9867450936054707
$ws.Range("Q2").Value = 0.034799143735
$ws.Range("R2").Value = 0.313192293615
$ws.Range("S2").Value = 0.0004809178690208855
$ws.Range("T2").Value = 0.0004823334993226473

# Row 3
$ws.Range("G3").Value = 0.101371
$ws.Range("H3").Value = 0.304113
$ws.Range("I3").Value = 0.0004873780190420389
$ws.Range("J3").Value = 0.0004888126654476159
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.004611333333333334
$ws.Range("N3").Value = 0.013834
$ws.Range("O3").Value = 0.0132549063945294
$ws.Range("P3").Value = 0.01325490639452941
$ws.Range("Q3").Value = 0.0004674554713333334
$ws.Range("R3").Value = 0.004207099242000001
$ws.Range("S3").Value = 0.000006460150021153395
$ws.Range("T3").Value = 0.000006479166124968568

# Row 4
$ws.Range("G4").Value = 123.540774
$ws.Range("H4").Value = 370.6223219999999
$ws.Range("I4").Value = 0.5939672855455723
$ws.Range("J4").Value = 0.5957156882185389
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.343285
$ws.Range("N4").Value = 1.029855
$ws.Range("O4").Value = 0.9867450936054706
$ws.Range("P4").Value = 0.9867450936054707
$ws.Range("Q4").Value = 42.40969460258999
$ws.Range("R4").Value = 381.6872514233099
$ws.Range("S4").Value = 0.5860943047742531
$ws.Range("T4").Value = 0.5878195325334495

# Row 5
$ws.Range("G5").Value = 123.540774
$ws.Range("H5").Value = 370.6223219999999
$ws.Range("I5").Value = 0.5939672855455723
$ws.Range("J5").Value = 0.5957156882185389
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.004611333333333334
$ws.Range("N5").Value = 0.013834
$ws.Range("O5").Value = 0.0132549063945294
$ws.Range("P5").Value = 0.01325490639452941
$ws.Range("Q5").Value = 0.569687689172
$ws.Range("R5").Value = 5.127189202547999
$ws.Range("S5").Value = 0.00787298077131928
$ws.Range("T5").Value = 0.007896155685089399

# Row 6
$ws.Range("G6").Value = 1.8313505
$ws.Range("H6").Value = 3.662701
$ws.Range("I6").Value = 0.008804884817764917
$ws.Range("J6").Value = 0.005887201923454927
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.343285
$ws.Range("N6").Value = 1.029855
$ws.Range("O6").Value = 0.9867450936054706
$ws.Range("P6").Value = 0.9867450936054707
$ws.Range("Q6").Value = 0.6286751563925
$ws.Range("R6").Value = 3.772050938355
$ws.Range("S6").Value = 0.00868817689369083
$ws.Range("T6").Value = 0.005809167613033839

# Row 7
$ws.Range("G7").Value = 1.8313505
$ws.Range("H7").Value = 3.662701
$ws.Range("I7").Value = 0.008804884817764917
$ws.Range("J7").Value = 0.005887201923454927
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.004611333333333334
$ws.Range("N7").Value = 0.013834
$ws.Range("O7").Value = 0.0132549063945294
$ws.Range("P7").Value = 0.01325490639452941
$ws.Range("Q7").Value = 0.008444967605666667
$ws.Range("R7").Value = 0.050669805634
$ws.Range("S7").Value = 0.0001167079240740871
$ws.Range("T7").Value = 0.00007803431042108854

# Row 8
$ws.Range("G8").Value = 82.51906066666666
$ws.Range("H8").Value = 247.557182
$ws.Range("I8").Value = 0.3967404516176207
$ws.Range("J8").Value = 0.3979082971925585
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.343285
$ws.Range("N8").Value = 1.029855
$ws.Range("O8").Value = 0.9867450936054706
$ws.Range("P8").Value = 0.9867450936054707
$ws.Range("Q8").Value = 28.32755574095667
$ws.Range("R8").Value = 254.94800166861
$ws.Range("S8").Value = 0.3914816940685058
$ws.Range("T8").Value = 0.3926340599596646

# Row 9
$ws.Range("G9").Value = 82.51906066666666
$ws.Range("H9").Value = 247.557182
$ws.Range("I9").Value = 0.3967404516176207
$ws.Range("J9").Value = 0.3979082971925585
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.004611333333333334
$ws.Range("N9").Value = 0.013834
$ws.Range("O9").Value = 0.0132549063945294
$ws.Range("P9").Value = 0.01325490639452941
$ws.Range("Q9").Value = 0.3805228950875555
$ws.Range("R9").Value = 3.424706055788
$ws.Range("S9").Value = 0.005258757549114884
$ws.Range("T9").Value = 0.005274237232893952

